$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix proveedor/name fields: commas used as separators were meant to be
#     periods (and "S.H." collapses to "SH" in one case). These are plain
#     text cells, so no numeric coercion to worry about - just overwrite
#     every cell that shared the old string with the corrected text.

$ws.Range("E43").Value  = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E166").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E177").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"

$ws.Range("E87").Value  = "FERNANDEZ. MARIO HUGO"

$ws.Range("E89").Value  = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

$ws.Range("E129").Value = "RICCOTTI. MARIANA EDITH"
$ws.Range("E142").Value = "RICCOTTI. MARIANA EDITH"

$ws.Range("E193").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

$ws.Range("F140").Value = "MERCANZINI. GASTON ARIEL"

# --- Fix "Importe" column: values were scraped with the Argentine
#     thousands/decimal punctuation ("1.234,56") still in the text instead
#     of being normalised to plain floating point ("1234.56"). Re-write
#     every amount in column H (rows 2-240) by stripping the thousands dot
#     and turning the decimal comma into a decimal point.
#
#     The column stores these amounts as literal text (not real numbers),
#     so the cells are temporarily switched to Text format while writing -
#     otherwise Excel would "helpfully" reinterpret "7170.00" as the number
#     7170 and drop the trailing zeros. The format is restored immediately
#     after so the workbook's styling is left as it was.

$importeRange = $ws.Range("H2:H240")
$importeRange.NumberFormat = "@"

for ($r = 2; $r -le 240; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $old = $cell.Text
    $new = $old.Replace(".", "").Replace(",", ".")
    $cell.Value = $new
}

$importeRange.Style = "Normal"
